$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.020.42'
$ws.Range("E2").Value = '  -1.92%  '

$ws.Range("D3").Value = '1.554.20'
$ws.Range("E3").Value = '  -1.20%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("E5").Value = '  +0.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '286.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3770'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3235'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.122'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.14'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -13.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07305'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.68%  '

$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.11%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.709'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.810'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.60%  '

$ws.Range("D16").Value = '1.549.60'
$ws.Range("E16").Value = '  -0.80%  '

$ws.Range("E17").Value = '  -3.30%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06636'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.31%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '84.91'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.429'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.60%  '

$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.86%  '

$ws.Range("D24").Value = '22.041.40'
$ws.Range("E24").Value = '  -1.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.237'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.18%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.518'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.30%  '

$ws.Range("E27").Value = '  -0.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.837'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.46%  '

$ws.Range("D30").Value = '1.726.18'
$ws.Range("E30").Value = '  -0.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.97%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.122'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.897'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08179'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.252'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.633'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -17.90%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.221'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02283'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06151'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2111'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.76%  '

$ws.Range("E41").Value = '  -7.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '10.88'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5932'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.96%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.62%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.722'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5728'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.929'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '119.73'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.155'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06905'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.02%  '
